# Automatische sync: 2025-06-17 14:57:43
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 15 ---
$ws.Cells.Item(15, 1).Value = "Re: Re: Re: Re: Wat zijn jullie openingstijden?"
$ws.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(15, 3).Value = @"
Geachte klant,
Bedankt voor uw interesse. Onze openingstijden zijn ma t/m vr van 9:00-18:00 en za van 10:00-16:00. Op zondag zijn we gesloten. Voor vragen zijn we bereikbaar via info@bedrijfsnaam.nl of telefonisch op [telefoonnummer].
Met vriendelijke groet,
[Naam van het bedrijf]
"@
$ws.Cells.Item(15, 4).Value = "Informatieaanvraag"
$ws.Cells.Item(15, 5).Value = @"
Geachte klant,
Dank voor uw bericht. Onze openingstijden zijn ma t/m vr van 9:00-18:00 en za van 10:00-16:00. Op zondag zijn wij gesloten. Voor vragen zijn wij bereikbaar via info@bedrijfsnaam.nl of telefonisch op [telefoonnummer].
Met vriendelijke groet,
[Naam van het bedrijf]
"@
$ws.Cells.Item(15, 6).Value = "2025-06-17 13:59:44"
$ws.Cells.Item(15, 7).Value = "Ja"

# --- Row 16 ---
$ws.Cells.Item(16, 1).Value = "Sollicitatie marketingfunctie"
$ws.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(16, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$ws.Cells.Item(16, 4).Value = "Overig"
$ws.Cells.Item(16, 6).Value = "2025-06-17 13:59:45"
$ws.Cells.Item(16, 7).Value = "Nee"

# --- Row 17 ---
$ws.Cells.Item(17, 1).Value = "Afmelding nieuwsbrief"
$ws.Cells.Item(17, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(17, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Cells.Item(17, 4).Value = "Afmelding"
$ws.Cells.Item(17, 6).Value = "2025-06-17 14:29:52"
$ws.Cells.Item(17, 7).Value = "Nee"

# --- Extend conditional formatting ranges to cover the new rows ---
$fcD = $ws.Range("D2:D14").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($ws.Range("D2:D17"))
$fcG = $ws.Range("G2:G14").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($ws.Range("G2:G17"))

# --- Update Dashboard summary counts ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 7
$dash.Cells.Item(3, 2).Value = 4
$dash.Cells.Item(4, 2).Value = 3
